# Generate Report for Handoff
# Adds two new file entries (477fe640-... and 862e4ef6-...) to the
# localization status report, on all three worksheets (Overview, zh-cn, de-de),
# inserted in order before the existing f48785cf-... entry.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

function Set-RowValues($ws, [int]$row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

function Add-StyledHyperlink($ws, [int]$row, [int]$col, $address, $display) {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, $col), $address, "", "", $display) | Out-Null
    $ws.Cells.Item($row, $col).Font.Underline = $true
    $ws.Cells.Item($row, $col).Font.Color = 0xED9564
}

function Set-DateFormat($ws, [int]$row, [int]$col) {
    $ws.Cells.Item($row, $col).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): columns A-G
#   A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#   E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------

Set-RowValues $wsOverview 4 @(
    "862e4ef6-4c76-491b-9069-ac64e2590c76.md",
    "e2e\862e4ef6-4c76-491b-9069-ac64e2590c76.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-08-24 10:42:33"
)

Set-RowValues $wsOverview 5 @(
    "f48785cf-9a83-4adb-a023-0e1f589f15d0.md",
    "e2e\f48785cf-9a83-4adb-a023-0e1f589f15d0.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-08-24 10:41:12"
)

# Row 3 (existing row) now becomes the 477fe640 entry.
Set-RowValues $wsOverview 3 @(
    "477fe640-84a0-475f-8064-ddd512e935ac.md",
    "e2e\477fe640-84a0-475f-8064-ddd512e935ac.md",
    ".md",
    "",
    "Ready for handoff",
    "Ready for handoff",
    "2016-08-24 10:42:33"
)

$wsOverview.Cells.Item(3,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item(5,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks on column B for rows 2-5 (clearing stale ones first).
$wsOverview.Hyperlinks.Delete()

Add-StyledHyperlink $wsOverview 2 2 "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/487bf27d8fa666a1088d35e77db8bdd279483b2a/e2e/836f7dee-297f-4940-88b2-4db5e13412a3.md" "e2e\836f7dee-297f-4940-88b2-4db5e13412a3.md"
Add-StyledHyperlink $wsOverview 3 2 "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e2537456d66c9b565a889dab2f9ddd0622a7373/e2e/477fe640-84a0-475f-8064-ddd512e935ac.md" "e2e\477fe640-84a0-475f-8064-ddd512e935ac.md"
Add-StyledHyperlink $wsOverview 4 2 "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eeb9356015983607ebf089a4efad77a5db6c60f0/e2e/862e4ef6-4c76-491b-9069-ac64e2590c76.md" "e2e\862e4ef6-4c76-491b-9069-ac64e2590c76.md"
Add-StyledHyperlink $wsOverview 5 2 "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/874fc95fabe2a71d0704b6fed0aafd34d6f435fa/e2e/f48785cf-9a83-4adb-a023-0e1f589f15d0.md" "e2e\f48785cf-9a83-4adb-a023-0e1f589f15d0.md"

# Resize the "Overview" table to include the new rows.
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) and "de-de" (sheet3): columns A-P
#   A=Source File Name, B=File Extension, C=Status, D=Source Path,
#   E=Priority, F=Content Duplicate, G=Latest Handoff File,
#   H=Latest Handoff Datetime, I=Latest Target File, J=Latest Handback File,
#   K=Latest Handback DateTime, L=Reference Tokens, M=To be localized,
#   N=Dependency From, O=Has metadata, P=Error Detail
# ---------------------------------------------------------------------------

Set-RowValues $wsZhCn 4 @(
    "862e4ef6-4c76-491b-9069-ac64e2590c76.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "862e4ef6-4c76-491b-9069-ac64e2590c76.824af31126bcc519d14bd95d05ee885a5b408bb8.zh-cn.xlf",
    "2016-08-24 10:42:28", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
Set-RowValues $wsZhCn 5 @(
    "f48785cf-9a83-4adb-a023-0e1f589f15d0.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "f48785cf-9a83-4adb-a023-0e1f589f15d0.8b7ade66b5f12130fb7b679359109bce73387788.zh-cn.xlf",
    "2016-08-24 10:41:02", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
# Row 3 becomes the 477fe640 entry.
Set-RowValues $wsZhCn 3 @(
    "477fe640-84a0-475f-8064-ddd512e935ac.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "477fe640-84a0-475f-8064-ddd512e935ac.9db682a3b17326c969623b587dd482ead88b4722.zh-cn.xlf",
    "2016-08-24 10:42:28", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)

Set-DateFormat $wsZhCn 3 8
Set-DateFormat $wsZhCn 3 11
Set-DateFormat $wsZhCn 4 8
Set-DateFormat $wsZhCn 4 11
Set-DateFormat $wsZhCn 5 8
Set-DateFormat $wsZhCn 5 11

$wsZhCn.Hyperlinks.Delete()

Add-StyledHyperlink $wsZhCn 2 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/487bf27d8fa666a1088d35e77db8bdd279483b2a/e2e/836f7dee-297f-4940-88b2-4db5e13412a3.md" "836f7dee-297f-4940-88b2-4db5e13412a3.md"
Add-StyledHyperlink $wsZhCn 2 9  "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5ee9ac9721378f4c01b82289b7d0ec9c83f8ef0b/e2e/836f7dee-297f-4940-88b2-4db5e13412a3.md" "836f7dee-297f-4940-88b2-4db5e13412a3.md"
Add-StyledHyperlink $wsZhCn 3 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e2537456d66c9b565a889dab2f9ddd0622a7373/e2e/477fe640-84a0-475f-8064-ddd512e935ac.md" "477fe640-84a0-475f-8064-ddd512e935ac.md"
Add-StyledHyperlink $wsZhCn 4 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eeb9356015983607ebf089a4efad77a5db6c60f0/e2e/862e4ef6-4c76-491b-9069-ac64e2590c76.md" "862e4ef6-4c76-491b-9069-ac64e2590c76.md"
Add-StyledHyperlink $wsZhCn 5 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/874fc95fabe2a71d0704b6fed0aafd34d6f435fa/e2e/f48785cf-9a83-4adb-a023-0e1f589f15d0.md" "f48785cf-9a83-4adb-a023-0e1f589f15d0.md"

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P5"))

Set-RowValues $wsDeDe 4 @(
    "862e4ef6-4c76-491b-9069-ac64e2590c76.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "862e4ef6-4c76-491b-9069-ac64e2590c76.824af31126bcc519d14bd95d05ee885a5b408bb8.de-de.xlf",
    "2016-08-24 10:42:33", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
Set-RowValues $wsDeDe 5 @(
    "f48785cf-9a83-4adb-a023-0e1f589f15d0.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "f48785cf-9a83-4adb-a023-0e1f589f15d0.8b7ade66b5f12130fb7b679359109bce73387788.de-de.xlf",
    "2016-08-24 10:41:12", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
# Row 3 becomes the 477fe640 entry.
Set-RowValues $wsDeDe 3 @(
    "477fe640-84a0-475f-8064-ddd512e935ac.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "477fe640-84a0-475f-8064-ddd512e935ac.9db682a3b17326c969623b587dd482ead88b4722.de-de.xlf",
    "2016-08-24 10:42:33", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)

Set-DateFormat $wsDeDe 3 8
Set-DateFormat $wsDeDe 3 11
Set-DateFormat $wsDeDe 4 8
Set-DateFormat $wsDeDe 4 11
Set-DateFormat $wsDeDe 5 8
Set-DateFormat $wsDeDe 5 11

$wsDeDe.Hyperlinks.Delete()

Add-StyledHyperlink $wsDeDe 2 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/487bf27d8fa666a1088d35e77db8bdd279483b2a/e2e/836f7dee-297f-4940-88b2-4db5e13412a3.md" "836f7dee-297f-4940-88b2-4db5e13412a3.md"
Add-StyledHyperlink $wsDeDe 2 9  "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b124c80ccfd3c66c74d6f3a4f4e4183eda7b3188/e2e/836f7dee-297f-4940-88b2-4db5e13412a3.md" "836f7dee-297f-4940-88b2-4db5e13412a3.md"
Add-StyledHyperlink $wsDeDe 3 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e2537456d66c9b565a889dab2f9ddd0622a7373/e2e/477fe640-84a0-475f-8064-ddd512e935ac.md" "477fe640-84a0-475f-8064-ddd512e935ac.md"
Add-StyledHyperlink $wsDeDe 4 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eeb9356015983607ebf089a4efad77a5db6c60f0/e2e/862e4ef6-4c76-491b-9069-ac64e2590c76.md" "862e4ef6-4c76-491b-9069-ac64e2590c76.md"
Add-StyledHyperlink $wsDeDe 5 1  "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/874fc95fabe2a71d0704b6fed0aafd34d6f435fa/e2e/f48785cf-9a83-4adb-a023-0e1f589f15d0.md" "f48785cf-9a83-4adb-a023-0e1f589f15d0.md"

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P5"))
